# Add "hydrogen combined cycle" as a new power plant type on the FSCaFoCC
# sheet, and rename the existing "hydrogen" row to "hydrogen combustion
# turbine" to disambiguate it from the new row.

$wb = $excel.ActiveWorkbook

$wsFS = $wb.Worksheets.Item("FSCaFoCC")

# Rename row 24's plant type label.
$wsFS.Range("A24").Value = "hydrogen combustion turbine"

# New row for the additional plant type.
$wsFS.Range("A25").Value = "hydrogen combined cycle"
$wsFS.Range("B25").Value = 0

# Match the formatting Excel applies to freshly-typed label cells here:
# an explicit black font color with vertically-centered alignment.
$label24 = $wsFS.Range("A24")
$label24.Font.Bold = $false
$label24.Font.Color = 0
$label24.VerticalAlignment = -4108

$label24.Copy()
$wsFS.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsFS.Range("E21").Select()

$wsFS.Activate()
